$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: years 2010-2014 with their corresponding values, replacing
# rows 2-16 (years 2000-2014) with just the last 5 years (rows 2-6)
$years = @("2010年", "2011年", "2012年", "2013年", "2014年")
$values = @(5293598.84491532, 1376574.77202485, 7977108.87673693, 6292500, 14013715.5670689)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Remove now-unused rows 7 through 16
$ws.Rows("7:16").Delete()
